{"js": "// Replace the 25 \"three-digit \u00d7 one-digit\" answer strings in the table\n// cells with their new values. Each old value is unique in the document,\n// so a body-wide literal search + replace (insertText with Replace)\n// for each pair is safe. Pairs are applied in document order so that a\n// later pair's new text never collides with an earlier pair's not-yet\n// processed old text.\nconst replacements = [\n  [\"425\u00d73=1275\", \"951\u00d78=7608\"],\n  [\"187\u00d77=1309\", \"357\u00d77=2499\"],\n  [\"127\u00d72=254\", \"786\u00d73=2358\"],\n  [\"833\u00d79=7497\", \"559\u00d77=3913\"],\n  [\"671\u00d73=2013\", \"188\u00d72=376\"],\n  [\"110\u00d79=990\", \"664\u00d75=3320\"],\n  [\"815\u00d75=4075\", \"487\u00d73=1461\"],\n  [\"507\u00d76=3042\", \"965\u00d75=4825\"],\n  [\"788\u00d75=3940\", \"103\u00d72=206\"],\n  [\"747\u00d78=5976\", \"129\u00d75=645\"],\n  [\"754\u00d77=5278\", \"175\u00d76=1050\"],\n  [\"641\u00d79=5769\", \"439\u00d76=2634\"],\n  [\"265\u00d74=1060\", \"960\u00d73=2880\"],\n  [\"627\u00d73=1881\", \"543\u00d74=2172\"],\n  [\"118\u00d72=236\", \"498\u00d76=2988\"],\n  [\"510\u00d79=4590\", \"114\u00d79=1026\"],\n  [\"405\u00d78=3240\", \"627\u00d73=1881\"],\n  [\"690\u00d74=2760\", \"498\u00d76=2988\"],\n  [\"739\u00d73=2217\", \"732\u00d73=2196\"],\n  [\"171\u00d74=684\", \"849\u00d72=1698\"],\n  [\"979\u00d78=7832\", \"639\u00d77=4473\"],\n  [\"504\u00d79=4536\", \"264\u00d72=528\"],\n  [\"347\u00d73=1041\", \"975\u00d75=4875\"],\n  [\"797\u00d73=2391\", \"578\u00d73=1734\"],\n  [\"115\u00d73=345\", \"658\u00d78=5264\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the 25 \"three-digit x one-digit\" answer strings in the table\n# cells with their new values. Each old value is unique in the document,\n# so Find/Replace (Execute with Replace:=wdReplaceOne, i.e. 2 -> replace\n# the single found occurrence) for each pair is safe. Pairs are applied\n# in document order so a later pair's new text never collides with an\n# earlier pair's not-yet-processed old text.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"425\u00d73=1275\", \"951\u00d78=7608\"),\n  @(\"187\u00d77=1309\", \"357\u00d77=2499\"),\n  @(\"127\u00d72=254\", \"786\u00d73=2358\"),\n  @(\"833\u00d79=7497\", \"559\u00d77=3913\"),\n  @(\"671\u00d73=2013\", \"188\u00d72=376\"),\n  @(\"110\u00d79=990\", \"664\u00d75=3320\"),\n  @(\"815\u00d75=4075\", \"487\u00d73=1461\"),\n  @(\"507\u00d76=3042\", \"965\u00d75=4825\"),\n  @(\"788\u00d75=3940\", \"103\u00d72=206\"),\n  @(\"747\u00d78=5976\", \"129\u00d75=645\"),\n  @(\"754\u00d77=5278\", \"175\u00d76=1050\"),\n  @(\"641\u00d79=5769\", \"439\u00d76=2634\"),\n  @(\"265\u00d74=1060\", \"960\u00d73=2880\"),\n  @(\"627\u00d73=1881\", \"543\u00d74=2172\"),\n  @(\"118\u00d72=236\", \"498\u00d76=2988\"),\n  @(\"510\u00d79=4590\", \"114\u00d79=1026\"),\n  @(\"405\u00d78=3240\", \"627\u00d73=1881\"),\n  @(\"690\u00d74=2760\", \"498\u00d76=2988\"),\n  @(\"739\u00d73=2217\", \"732\u00d73=2196\"),\n  @(\"171\u00d74=684\", \"849\u00d72=1698\"),\n  @(\"979\u00d78=7832\", \"639\u00d77=4473\"),\n  @(\"504\u00d79=4536\", \"264\u00d72=528\"),\n  @(\"347\u00d73=1041\", \"975\u00d75=4875\"),\n  @(\"797\u00d73=2391\", \"578\u00d73=1734\"),\n  @(\"115\u00d73=345\", \"658\u00d78=5264\")\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
